$wb = $excel.ActiveWorkbook

# Keep a handle on a sheet that already has the "bold/border/centered" header
# style (style id 2) and the "bold/border" index style used on column A data
# cells, so we can copy that formatting onto newly written cells instead of
# guessing at font/border/alignment settings by hand.
$fmtSrc = $wb.Worksheets.Item("2021-Q2")

# --- Step 1: the existing "总计" (grand-totals) sheet is repurposed in place
#     to become the new "2022-Q1" per-fund-holdings sheet (same sheetId/rId). ---
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Range("A1:H10").Clear()

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$fmtSrc.Range("B1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$q1.Range("A2").Value = 0
$fmtSrc.Range("A2").Copy()
$q1.Range("A2").PasteSpecial(-4122)

$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "005295"
$q1.Range("B2").Style = "Normal"

$q1.Range("C2").Value = "诺德天富灵活配置混合"

$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value = "1.21"
$q1.Range("D2").Style = "Normal"

$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value = "93.81"
$q1.Range("E2").Style = "Normal"

$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value = "1.28"
$q1.Range("F2").Style = "Normal"

$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value = "0.0155"
$q1.Range("G2").Style = "Normal"

$q1.Range("H2").Value = 7

# --- Step 2: insert a brand-new "总计" sheet right after "2022-Q1", rebuilt
#     from a copy of an existing sheet so it keeps the same sheetPr/margins
#     boilerplate, then wiped and refilled with the updated totals table. ---
$fmtSrc.Copy([System.Reflection.Missing]::Value, $q1)
$total = $wb.Worksheets.Item($q1.Index + 1)
$total.Name = "总计"
$total.Range("A1:H10").Clear()

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$fmtSrc.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.02

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q2"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.12

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q1"
$total.Range("C4").Value = 3
$total.Range("D4").Value = 0.4

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2020-Q4"
$total.Range("C5").Value = 3
$total.Range("D5").Value = 0.38

$fmtSrc.Range("A2").Copy()
$total.Range("A2:A5").PasteSpecial(-4122)
